$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 ("I0") and J1 ("IF"), styled to match the existing
# header cells (bold font, thin border all around, centered horizontally,
# top-aligned vertically) like B1:H1.
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$header = $ws.Range("I1:J1")
$header.Font.Bold = $true
$header.Borders.LineStyle = 1
$header.HorizontalAlignment = -4108
$header.VerticalAlignment = -4160

# Data values for the new columns I (I0) and J (IF), rows 2-19.
$data = @(
    @(8, 8),
    @(5, 5),
    @(4, 6),
    @(9, 9),
    @(8, 8),
    @(10, 10),
    @(8, 8),
    @(6, 7),
    @(9, 9),
    @(6, 6),
    @(2, 2),
    @(7, 7),
    @(8, 9),
    @(9, 9),
    @(5, 5),
    @(7, 7),
    @(6, 6),
    @(6, 6)
)

$row = 2
foreach ($pair in $data) {
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
    $row++
}
